# "User testing.docx" edit: turn the feedback items into a numbered list,
# drop three stale/duplicate notes, fix up a couple of sentences, and fold
# in some additional feedback text.

$d = $word.ActiveDocument

# --- 1. Remove paragraphs that were dropped in the revision -----------------
# (delete from the bottom up so earlier paragraph indices stay valid)

# "We need to fill the white space on the front page right now with something
#  spectacular (since there's room)."
$d.Paragraphs(11).Range.Delete()

# "When you hit a button to move a show (or add a show from a search) you can
#  tell by the count but the button gives no indication it was pressed. Was
#  confusing."
$d.Paragraphs(5).Range.Delete()

# "Content management -if you're on "want to watch" then there should only be
#  buttons for Watched or Watching, omit the button for the tab you're on"
$d.Paragraphs(3).Range.Delete()

# After the deletions above, the surviving feedback paragraphs are now
# contiguous: 2 Inconsistent Spanish switch, 3 Binge number, 4 Mardi Gras,
# 5 Confirm notifications, 6 For shows currently playing, 7 If we prompt for
# a name, 8 A "Pro" toggle.

# --- 2. Append new sentence to the "Binge number" paragraph -----------------
$p = $d.Paragraphs(3)
$r = $p.Range
$ins = $d.Range($r.End - 1, $r.End - 1)
$ins.InsertAfter(". Let’s remove it from the front page for now. ")

# --- 3. Fix up the "Mardi Gras" paragraph wording ----------------------------
$p = $d.Paragraphs(4)
$r = $p.Range
$r.Find.Execute("which is  a good look", $false, $false, $false, $false, $false, $true, 1, $false, "which is a good look", 2)
$r2 = $p.Range
$r2.Find.Execute("the buttons those", $false, $false, $false, $false, $false, $true, 1, $false, "the buttons with those", 2)

# --- 4. Append new sentences to the "If we prompt for a name" paragraph -----
$p = $d.Paragraphs(7)
$r = $p.Range
$ins = $d.Range($r.End - 1, $r.End - 1)
$ins.InsertAfter("So if it’s Travis right now, but someone in settings types in “Fredrico” there should be a prompt that says “This is Frederico’s now, the only way back is to come back here.” Or something similarly snarky. ")

# --- 5. Turn the feedback paragraphs (2-8) into a numbered list -------------
$start = $d.Paragraphs(2).Range.Start
$end = $d.Paragraphs(8).Range.End
$listRange = $d.Range($start, $end)
$listRange.Style = "List Paragraph"
$listRange.ListFormat.ApplyNumberDefault()
